# edit.ps1
#
# Reproduces (the COM-reachable parts of) the target diff:
#   - WithTable            : D2:D5 boolean cells drop their local date-style
#                             override (s="1" -> default style) and the
#                             sheet selection moves to E23.
#   - Tableless             : D2:D5 boolean cells drop their local date-style
#                             override (s="1" -> default style), the sheet
#                             selection moves to I19, and it stops being the
#                             active/tab-selected sheet.
#   - WithTable_Duplicate   : E5:E8 boolean cells drop their local date-style
#                             override (s="1" -> default style), the sheet
#                             selection moves to M11, and it becomes the
#                             active/tab-selected sheet (so workbook-level
#                             activeTab moves from the 2nd to the 3rd tab).
#
# NOTE: a few items in the source diff correspond to metadata that isn't
# exposed anywhere in the Excel object model surfaced by this host (window
# pixel dimensions / xr:revisionPtr document id / the x15ac absolute-path
# breadcrumb / the cell-style "Normal"->"Standard" rename / the dxf-only
# table-column bookkeeping / theme display name) - those are artifacts of
# the file being round-tripped through a different authoring tool, not
# content an Excel user/script can set, so they're intentionally left
# alone here.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("WithTable")
$ws2 = $wb.Worksheets.Item("Tableless")
$ws3 = $wb.Worksheets.Item("WithTable_Duplicate")

# Reset the "Boolean" column cells back to the default cell style (removes
# the inherited date number-format override that was left on them).
$ws1.Range("D2:D5").Style = "Normal"
$ws2.Range("D2:D5").Style = "Normal"
$ws3.Range("E5:E8").Style = "Normal"

# Update each sheet's remembered selection. Selecting a range also makes
# that sheet the active one, so WithTable_Duplicate (the sheet that should
# end up active/tab-selected) is selected last.
$ws1.Range("E23").Select()
$ws2.Range("I19").Select()
$ws3.Range("M11").Select()
